# Fruta / hortaliza, semanal
#
# Insert a new weekly price record as row 47 (pushing the existing rows
# 47-53 down to 48-54) on the single worksheet of the workbook.
# The new row duplicates the commercial data of the (old) row 48 record
# but carries a newer date (2022-08-03 -> Excel serial 44776).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 47:53 down to 48:54, leaving a blank row 47 in their place.
$ws.Rows("47:47").Insert()

# Populate the newly inserted row 47.
$ws.Range("A47").Value = 9
$ws.Range("B47").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C47").Value = "Metropolitana"
$ws.Range("D47").Value = 44776
$ws.Range("E47").Value = 13
$ws.Range("F47").Value = 100112035
$ws.Range("G47").Value = "Bruselas (repollito)"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 34
$ws.Range("K47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("M47").Value = 20000
$ws.Range("N47").Value = "$/malla 15 kilos"
$ws.Range("O47").Value = "Hijuelas"
$ws.Range("P47").Value = 1333
$ws.Range("Q47").Value = 15
$ws.Range("R47").Value = "Hortaliza"
